# Add 5 new benchmark rows (tree-model test runs on the avg MoCo tile
# features) to the "Tabelle1" table, expanding it from A1:O31 to A1:O36,
# and update the sheet view to reflect the new selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the columns that reuse already-existing shared strings / are
#     plain numbers first (order among these does not affect the shared
#     string table). ---

# Row 32 - XGB Classifier
$ws.Range("B32").Value = "XGB Classifier"
$ws.Range("C32").Value = "MoCo"
$ws.Range("D32").Value = "Centers"
$ws.Range("E32").Value = "1 x 3"
$ws.Range("F32").Value = "average"
$ws.Range("J32").Value = 0.654
$ws.Range("K32").Value = 0.696
$ws.Range("L32").Value = 0.599

# Row 33 - Catboost
$ws.Range("C33").Value = "MoCo"
$ws.Range("D33").Value = "Centers"
$ws.Range("E33").Value = "1 x 3"
$ws.Range("F33").Value = "average"
$ws.Range("J33").Value = 0.63
$ws.Range("K33").Value = 0.704
$ws.Range("L33").Value = 0.597

# Row 34 - LightGBM
$ws.Range("C34").Value = "MoCo"
$ws.Range("D34").Value = "Centers"
$ws.Range("E34").Value = "1 x 3"
$ws.Range("F34").Value = "average"
$ws.Range("J34").Value = 0.669
$ws.Range("K34").Value = 0.698
$ws.Range("L34").Value = 0.595

# Row 35 - ExtraTrees
$ws.Range("C35").Value = "MoCo"
$ws.Range("D35").Value = "Centers"
$ws.Range("E35").Value = "1 x 3"
$ws.Range("F35").Value = "average"
$ws.Range("J35").Value = 0.631
$ws.Range("K35").Value = 0.645
$ws.Range("L35").Value = 0.595

# Row 36 - DecisionTree
$ws.Range("C36").Value = "MoCo"
$ws.Range("D36").Value = "Centers"
$ws.Range("E36").Value = "1 x 3"
$ws.Range("F36").Value = "average"
$ws.Range("J36").Value = 0.505
$ws.Range("K36").Value = 0.563
$ws.Range("L36").Value = 0.529

# --- Now write the "Model" (B) and "Average Val AUC" (M) cells, the only
#     ones that introduce brand-new shared strings, in the precise order
#     the values were originally entered. ---

$ws.Range("M32").Value = "0.649 (0.040)"
$ws.Range("M33").Value = "0.644 (0.045)"
$ws.Range("B33").Value = "Catboost"
$ws.Range("M34").Value = "0.654 (0.043)"
$ws.Range("B34").Value = "LightGBM"
$ws.Range("B35").Value = "ExtraTrees"
$ws.Range("M35").Value = "0.624 (0.021)"
$ws.Range("B36").Value = "DecisionTree"
$ws.Range("M36").Value = "0.532 (0.024)"

# --- Expand the table / autofilter range to cover the new rows. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O36"))

# --- Update the sheet view: scrolled position + active selection. ---
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("N36").Select()

Write-Host "done"
